$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 159; existing rows 159..189 shift down to 160..190.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with a new data record (same constant
# columns as the surrounding block, new D/J/K/L/M/P values).
$ws.Cells.Item(159, 1).Value = 8
$ws.Cells.Item(159, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(159, 3).Value = "Coquimbo"
$ws.Cells.Item(159, 4).Value = 44694
$ws.Cells.Item(159, 5).Value = 4
$ws.Cells.Item(159, 6).Value = 100112037
$ws.Cells.Item(159, 7).Value = "Cebollín"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 2300
$ws.Cells.Item(159, 11).Value = 1100
$ws.Cells.Item(159, 12).Value = 1200
$ws.Cells.Item(159, 13).Value = 1150
$ws.Cells.Item(159, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(159, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(159, 16).Value = 192
$ws.Cells.Item(159, 17).Value = 6
$ws.Cells.Item(159, 18).Value = "Hortaliza"
